# Replace each arithmetic "old" expression in the answer table with its
# corresponding "new" expression, exactly once and in document order.
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
$d = $word.ActiveDocument
$d.Content.Find.Execute("70-64=6", $true, $true, $false, $false, $false, $true, 1, $false, "77-72=5", 2) | Out-Null
$d.Content.Find.Execute("3+63=66", $true, $true, $false, $false, $false, $true, 1, $false, "68-20=48", 2) | Out-Null
$d.Content.Find.Execute("4+52=56", $true, $true, $false, $false, $false, $true, 1, $false, "50+1=51", 2) | Out-Null
$d.Content.Find.Execute("59-7=52", $true, $true, $false, $false, $false, $true, 1, $false, "31+44=75", 2) | Out-Null
$d.Content.Find.Execute("7+55=62", $true, $true, $false, $false, $false, $true, 1, $false, "87-22=65", 2) | Out-Null
$d.Content.Find.Execute("36+21=57", $true, $true, $false, $false, $false, $true, 1, $false, "25+56=81", 2) | Out-Null
$d.Content.Find.Execute("43-5=38", $true, $true, $false, $false, $false, $true, 1, $false, "57-7=50", 2) | Out-Null
$d.Content.Find.Execute("26+70=96", $true, $true, $false, $false, $false, $true, 1, $false, "83-18=65", 2) | Out-Null
$d.Content.Find.Execute("93-58=35", $true, $true, $false, $false, $false, $true, 1, $false, "27+71=98", 2) | Out-Null
$d.Content.Find.Execute("88-81=7", $true, $true, $false, $false, $false, $true, 1, $false, "85-84=1", 2) | Out-Null
$d.Content.Find.Execute("38+35=73", $true, $true, $false, $false, $false, $true, 1, $false, "68-26=42", 2) | Out-Null
$d.Content.Find.Execute("1+76=77", $true, $true, $false, $false, $false, $true, 1, $false, "97-2=95", 2) | Out-Null
$d.Content.Find.Execute("19+2=21", $true, $true, $false, $false, $false, $true, 1, $false, "57-46=11", 2) | Out-Null
$d.Content.Find.Execute("2+7=9", $true, $true, $false, $false, $false, $true, 1, $false, "54-0=54", 2) | Out-Null
$d.Content.Find.Execute("6+55=61", $true, $true, $false, $false, $false, $true, 1, $false, "88-5=83", 2) | Out-Null
$d.Content.Find.Execute("87-77=10", $true, $true, $false, $false, $false, $true, 1, $false, "28+56=84", 2) | Out-Null
$d.Content.Find.Execute("77+13=90", $true, $true, $false, $false, $false, $true, 1, $false, "20-3=17", 2) | Out-Null
$d.Content.Find.Execute("48+36=84", $true, $true, $false, $false, $false, $true, 1, $false, "74+0=74", 2) | Out-Null
$d.Content.Find.Execute("14+57=71", $true, $true, $false, $false, $false, $true, 1, $false, "90-83=7", 2) | Out-Null
$d.Content.Find.Execute("58+41=99", $true, $true, $false, $false, $false, $true, 1, $false, "87-80=7", 2) | Out-Null
$d.Content.Find.Execute("21+13=34", $true, $true, $false, $false, $false, $true, 1, $false, "23+69=92", 2) | Out-Null
$d.Content.Find.Execute("99-59=40", $true, $true, $false, $false, $false, $true, 1, $false, "39+31=70", 2) | Out-Null
$d.Content.Find.Execute("36-1=35", $true, $true, $false, $false, $false, $true, 1, $false, "81-52=29", 2) | Out-Null
$d.Content.Find.Execute("46-6=40", $true, $true, $false, $false, $false, $true, 1, $false, "81-75=6", 2) | Out-Null
$d.Content.Find.Execute("47+43=90", $true, $true, $false, $false, $false, $true, 1, $false, "35+42=77", 2) | Out-Null
$d.Content.Find.Execute("94-54=40", $true, $true, $false, $false, $false, $true, 1, $false, "87-80=7", 2) | Out-Null
$d.Content.Find.Execute("12+8=20", $true, $true, $false, $false, $false, $true, 1, $false, "36+7=43", 2) | Out-Null
$d.Content.Find.Execute("8+71=79", $true, $true, $false, $false, $false, $true, 1, $false, "37+34=71", 2) | Out-Null
$d.Content.Find.Execute("63+12=75", $true, $true, $false, $false, $false, $true, 1, $false, "23+36=59", 2) | Out-Null
$d.Content.Find.Execute("87-39=48", $true, $true, $false, $false, $false, $true, 1, $false, "34+44=78", 2) | Out-Null
$d.Content.Find.Execute("15+21=36", $true, $true, $false, $false, $false, $true, 1, $false, "21-18=3", 2) | Out-Null
$d.Content.Find.Execute("25+61=86", $true, $true, $false, $false, $false, $true, 1, $false, "74-61=13", 2) | Out-Null
$d.Content.Find.Execute("89-22=67", $true, $true, $false, $false, $false, $true, 1, $false, "58+3=61", 2) | Out-Null
$d.Content.Find.Execute("6+80=86", $true, $true, $false, $false, $false, $true, 1, $false, "71-51=20", 2) | Out-Null
$d.Content.Find.Execute("37-2=35", $true, $true, $false, $false, $false, $true, 1, $false, "94-64=30", 2) | Out-Null
$d.Content.Find.Execute("88-43=45", $true, $true, $false, $false, $false, $true, 1, $false, "33+50=83", 2) | Out-Null
$d.Content.Find.Execute("23-0=23", $true, $true, $false, $false, $false, $true, 1, $false, "45+35=80", 2) | Out-Null
$d.Content.Find.Execute("26-1=25", $true, $true, $false, $false, $false, $true, 1, $false, "57+21=78", 2) | Out-Null
$d.Content.Find.Execute("15+64=79", $true, $true, $false, $false, $false, $true, 1, $false, "87-31=56", 2) | Out-Null
$d.Content.Find.Execute("44+29=73", $true, $true, $false, $false, $false, $true, 1, $false, "63+16=79", 2) | Out-Null
$d.Content.Find.Execute("29+68=97", $true, $true, $false, $false, $false, $true, 1, $false, "75+16=91", 2) | Out-Null
$d.Content.Find.Execute("8+42=50", $true, $true, $false, $false, $false, $true, 1, $false, "69-51=18", 2) | Out-Null
$d.Content.Find.Execute("9+75=84", $true, $true, $false, $false, $false, $true, 1, $false, "61-9=52", 2) | Out-Null
$d.Content.Find.Execute("22+12=34", $true, $true, $false, $false, $false, $true, 1, $false, "28-6=22", 2) | Out-Null
$d.Content.Find.Execute("99-55=44", $true, $true, $false, $false, $false, $true, 1, $false, "18+32=50", 2) | Out-Null
$d.Content.Find.Execute("4+45=49", $true, $true, $false, $false, $false, $true, 1, $false, "80-37=43", 2) | Out-Null
$d.Content.Find.Execute("51+33=84", $true, $true, $false, $false, $false, $true, 1, $false, "66+2=68", 2) | Out-Null
$d.Content.Find.Execute("44+40=84", $true, $true, $false, $false, $false, $true, 1, $false, "97-26=71", 2) | Out-Null
$d.Content.Find.Execute("75-13=62", $true, $true, $false, $false, $false, $true, 1, $false, "11+28=39", 2) | Out-Null
$d.Content.Find.Execute("18+28=46", $true, $true, $false, $false, $false, $true, 1, $false, "75-62=13", 2) | Out-Null
$d.Content.Find.Execute("19-18=1", $true, $true, $false, $false, $false, $true, 1, $false, "93-29=64", 2) | Out-Null
$d.Content.Find.Execute("68-30=38", $true, $true, $false, $false, $false, $true, 1, $false, "29+59=88", 2) | Out-Null
$d.Content.Find.Execute("41-22=19", $true, $true, $false, $false, $false, $true, 1, $false, "57+38=95", 2) | Out-Null
$d.Content.Find.Execute("60-16=44", $true, $true, $false, $false, $false, $true, 1, $false, "87+5=92", 2) | Out-Null
$d.Content.Find.Execute("43+56=99", $true, $true, $false, $false, $false, $true, 1, $false, "49+44=93", 2) | Out-Null
$d.Content.Find.Execute("58-52=6", $true, $true, $false, $false, $false, $true, 1, $false, "3+1=4", 2) | Out-Null
$d.Content.Find.Execute("53+34=87", $true, $true, $false, $false, $false, $true, 1, $false, "91-8=83", 2) | Out-Null
$d.Content.Find.Execute("87-15=72", $true, $true, $false, $false, $false, $true, 1, $false, "57-0=57", 2) | Out-Null
$d.Content.Find.Execute("79+16=95", $true, $true, $false, $false, $false, $true, 1, $false, "66-49=17", 2) | Out-Null
$d.Content.Find.Execute("64+0=64", $true, $true, $false, $false, $false, $true, 1, $false, "80-41=39", 2) | Out-Null
$d.Content.Find.Execute("82-75=7", $true, $true, $false, $false, $false, $true, 1, $false, "68-4=64", 2) | Out-Null
$d.Content.Find.Execute("31+12=43", $true, $true, $false, $false, $false, $true, 1, $false, "50-44=6", 2) | Out-Null
$d.Content.Find.Execute("33+22=55", $true, $true, $false, $false, $false, $true, 1, $false, "16+34=50", 2) | Out-Null
$d.Content.Find.Execute("28+42=70", $true, $true, $false, $false, $false, $true, 1, $false, "84+4=88", 2) | Out-Null
$d.Content.Find.Execute("18+55=73", $true, $true, $false, $false, $false, $true, 1, $false, "15+80=95", 2) | Out-Null
$d.Content.Find.Execute("95-21=74", $true, $true, $false, $false, $false, $true, 1, $false, "68-23=45", 2) | Out-Null
$d.Content.Find.Execute("12+12=24", $true, $true, $false, $false, $false, $true, 1, $false, "49+26=75", 2) | Out-Null
$d.Content.Find.Execute("20-1=19", $true, $true, $false, $false, $false, $true, 1, $false, "69-35=34", 2) | Out-Null
$d.Content.Find.Execute("55+7=62", $true, $true, $false, $false, $false, $true, 1, $false, "47+4=51", 2) | Out-Null
$d.Content.Find.Execute("89-16=73", $true, $true, $false, $false, $false, $true, 1, $false, "3+45=48", 2) | Out-Null
$d.Content.Find.Execute("82-44=38", $true, $true, $false, $false, $false, $true, 1, $false, "93-31=62", 2) | Out-Null
$d.Content.Find.Execute("16-7=9", $true, $true, $false, $false, $false, $true, 1, $false, "35-26=9", 2) | Out-Null
$d.Content.Find.Execute("79-22=57", $true, $true, $false, $false, $false, $true, 1, $false, "58-49=9", 2) | Out-Null
$d.Content.Find.Execute("88-15=73", $true, $true, $false, $false, $false, $true, 1, $false, "23+36=59", 2) | Out-Null
$d.Content.Find.Execute("78-66=12", $true, $true, $false, $false, $false, $true, 1, $false, "52-41=11", 2) | Out-Null
$d.Content.Find.Execute("17-5=12", $true, $true, $false, $false, $false, $true, 1, $false, "43-16=27", 2) | Out-Null
$d.Content.Find.Execute("87-40=47", $true, $true, $false, $false, $false, $true, 1, $false, "26-6=20", 2) | Out-Null
$d.Content.Find.Execute("26+32=58", $true, $true, $false, $false, $false, $true, 1, $false, "64-10=54", 2) | Out-Null
$d.Content.Find.Execute("32-26=6", $true, $true, $false, $false, $false, $true, 1, $false, "55-35=20", 2) | Out-Null
$d.Content.Find.Execute("4+21=25", $true, $true, $false, $false, $false, $true, 1, $false, "74-72=2", 2) | Out-Null
$d.Content.Find.Execute("54-6=48", $true, $true, $false, $false, $false, $true, 1, $false, "47+17=64", 2) | Out-Null
$d.Content.Find.Execute("54+13=67", $true, $true, $false, $false, $false, $true, 1, $false, "2+46=48", 2) | Out-Null
$d.Content.Find.Execute("96-24=72", $true, $true, $false, $false, $false, $true, 1, $false, "33+6=39", 2) | Out-Null
$d.Content.Find.Execute("62+30=92", $true, $true, $false, $false, $false, $true, 1, $false, "4+47=51", 2) | Out-Null
$d.Content.Find.Execute("74-37=37", $true, $true, $false, $false, $false, $true, 1, $false, "29+37=66", 2) | Out-Null
$d.Content.Find.Execute("99-91=8", $true, $true, $false, $false, $false, $true, 1, $false, "59-11=48", 2) | Out-Null
$d.Content.Find.Execute("52-31=21", $true, $true, $false, $false, $false, $true, 1, $false, "69+27=96", 2) | Out-Null
$d.Content.Find.Execute("0+67=67", $true, $true, $false, $false, $false, $true, 1, $false, "40+9=49", 2) | Out-Null
$d.Content.Find.Execute("37+27=64", $true, $true, $false, $false, $false, $true, 1, $false, "34+56=90", 2) | Out-Null
$d.Content.Find.Execute("79-25=54", $true, $true, $false, $false, $false, $true, 1, $false, "8+24=32", 2) | Out-Null
$d.Content.Find.Execute("87-35=52", $true, $true, $false, $false, $false, $true, 1, $false, "38+56=94", 2) | Out-Null
$d.Content.Find.Execute("80+4=84", $true, $true, $false, $false, $false, $true, 1, $false, "8-4=4", 2) | Out-Null
$d.Content.Find.Execute("62-7=55", $true, $true, $false, $false, $false, $true, 1, $false, "66-6=60", 2) | Out-Null
$d.Content.Find.Execute("36+26=62", $true, $true, $false, $false, $false, $true, 1, $false, "84-64=20", 2) | Out-Null
$d.Content.Find.Execute("69-23=46", $true, $true, $false, $false, $false, $true, 1, $false, "24-7=17", 2) | Out-Null
$d.Content.Find.Execute("70-61=9", $true, $true, $false, $false, $false, $true, 1, $false, "7+10=17", 2) | Out-Null
$d.Content.Find.Execute("30+68=98", $true, $true, $false, $false, $false, $true, 1, $false, "39-24=15", 2) | Out-Null
$d.Content.Find.Execute("96-10=86", $true, $true, $false, $false, $false, $true, 1, $false, "48-22=26", 2) | Out-Null
$d.Content.Find.Execute("84-46=38", $true, $true, $false, $false, $false, $true, 1, $false, "40+14=54", 2) | Out-Null
$d.Content.Find.Execute("64+19=83", $true, $true, $false, $false, $false, $true, 1, $false, "69-1=68", 2) | Out-Null
